$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.230.64"
$ws.Range("E2").Value = "  +0.41%  "

# Row 3
$ws.Range("D3").Value = "1.857.90"
$ws.Range("E3").Value = "  +0.36%  "

# Row 4
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").Value = "'0.7098"
$ws.Range("E5").Value = "  +2.06%  "

# Row 6
$ws.Range("D6").Value = "'238.63"
$ws.Range("E6").Value = "  +0.10%  "

# Row 7
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("D8").Value = "'0.08002"
$ws.Range("E8").Value = "  +3.64%  "

# Row 9
$ws.Range("D9").Value = "'0.3033"
$ws.Range("E9").Value = "  -0.12%  "

# Row 10
$ws.Range("D10").Value = "'23.49"
$ws.Range("E10").Value = "  +0.49%  "

# Row 11
$ws.Range("D11").Value = "'0.08203"
$ws.Range("E11").Value = "  +1.05%  "

# Row 12
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'5.176"
$ws.Range("E12").Value = "  -0.72%  "

# Row 13
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'0.7043"
$ws.Range("E13").Value = "  -3.18%  "

# Row 14
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").Value = "'89.61"
$ws.Range("E14").Value = "  +0.63%  "

# Row 15
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.689.18"
$ws.Range("E15").Value = "  -9.01%  "

# Row 16
$ws.Range("D16").Value = "29.186.34"
$ws.Range("E16").Value = "  +0.31%  "

# Row 17
$ws.Range("D17").Value = "'5.827"
$ws.Range("E17").Value = "  +1.12%  "

# Row 18
$ws.Range("D18").Value = "'0.000007871"
$ws.Range("E18").Value = "  +1.76%  "

# Row 19
$ws.Range("D19").Value = "'13.26"
$ws.Range("E19").Value = "  +0.15%  "

# Row 20
$ws.Range("D20").Value = "'238.04"
$ws.Range("E20").Value = "  +0.40%  "

# Row 21
$ws.Range("D21").Value = "'0.9984"
$ws.Range("E21").Value = "  -0.08%  "

# Row 22
$ws.Range("D22").Value = "2.082.95"
$ws.Range("E22").Value = "  -0.52%  "

# Row 23
$ws.Range("E23").Value = "  +0.06%  "

# Row 24
$ws.Range("D24").Value = "'7.493"
$ws.Range("E24").Value = "  -1.39%  "

# Row 25
$ws.Range("D25").Value = "'162.91"
$ws.Range("E25").Value = "  +1.10%  "

# Row 26
$ws.Range("D26").Value = "'8.925"
$ws.Range("E26").Value = "  -0.66%  "

# Row 27
$ws.Range("D27").Value = "'0.1445"
$ws.Range("E27").Value = "  +0.47%  "

# Row 28
$ws.Range("E28").Value = "  +0.45%  "

# Row 29
$ws.Range("D29").Value = "'1.919"
$ws.Range("E29").Value = "  -3.59%  "

# Row 30
$ws.Range("D30").Value = "'1.428"
$ws.Range("E30").Value = "  +1.68%  "

# Row 31
$ws.Range("E31").Value = "  -0.45%  "

# Row 32
$ws.Range("E32").Value = "  -2.88%  "

# Row 33
$ws.Range("E33").Value = "  -0.01%  "

# Row 34
$ws.Range("D34").Value = "'0.05205"
$ws.Range("E34").Value = "  -0.68%  "

# Row 35
$ws.Range("D35").Value = "'1.160"
$ws.Range("E35").Value = "  -2.36%  "

# Row 36
$ws.Range("D36").Value = "'0.7112"
$ws.Range("E36").Value = "  +1.19%  "

# Row 37
$ws.Range("D37").Value = "'0.9994"
$ws.Range("E37").Value = "  -2.49%  "

# Row 38
$ws.Range("D38").Value = "'2.675"
$ws.Range("E38").Value = "  +1.05%  "

# Row 39
$ws.Range("E39").Value = "  -0.06%  "

# Row 40
$ws.Range("D40").Value = "'2.720"
$ws.Range("E40").Value = "  +1.64%  "

# Row 41
$ws.Range("D41").Value = "'0.9302"
$ws.Range("E41").Value = "  +0.65%  "

# Row 42
$ws.Range("D42").Value = "1.135.12"
$ws.Range("E42").Value = "  +5.02%  "

# Row 43
$ws.Range("D43").Value = "'0.4256"
$ws.Range("E43").Value = "  -0.28%  "

# Row 44
$ws.Range("D44").Value = "'5.885"
$ws.Range("E44").Value = "  -1.81%  "

# Row 45
$ws.Range("D45").Value = "'70.04"
$ws.Range("E45").Value = "  -0.52%  "

# Row 46
$ws.Range("D46").Value = "'0.9999"
$ws.Range("E46").Value = "  +0.00%  "

# Row 47
$ws.Range("D47").Value = "'102.79"
$ws.Range("E47").Value = "  -0.41%  "

# Row 48
$ws.Range("D48").Value = "'0.5336"
$ws.Range("E48").Value = "  -4.28%  "

# Row 49
$ws.Range("D49").Value = "'1.768"
$ws.Range("E49").Value = "  -0.53%  "

# Row 50
$ws.Range("D50").Value = "1.979.40"
$ws.Range("E50").Value = "  -0.56%  "

# Row 51
$ws.Range("D51").Value = "'9.163"
$ws.Range("E51").Value = "  +0.00%  "
